$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.993.64"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.556.49"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.01"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.248"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.55"
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0860"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.777.56"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.559.32"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.70"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.515"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.984.33"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.68"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.57"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0687"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.25"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("E24").Value = "  -2.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.34"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.88"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0460"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.371.13"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.969"
$ws.Range("E36").Value = "  +5.67%  "
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.818"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.521"
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.979"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.48"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.24"
$ws.Range("E44").Value = "  +2.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.95"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.75"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.690.37"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.24"
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.30"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0958"
$ws.Range("E51").Value = "  +0.61%  "
